$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 33-35: Results column (E) changes from PASS to SKIP
$ws.Range("E33").Value = "SKIP"
$ws.Range("E34").Value = "SKIP"
$ws.Range("E35").Value = "SKIP"

# New row 36: copy formatting from row 35, then set values
$ws.Range("A35:E35").Copy()
$ws.Range("A36:E36").PasteSpecial(-4122)
$ws.Range("A36").Value = "PublishedAPostLikeCountTest"
$ws.Range("B36").Value = "TBD"
$ws.Range("C36").Value = "Verify that POST tab count getting increased while appreciate post from Record view page"
$ws.Range("D36").Value = "Y"
$ws.Range("E36").Value = "SKIP"

# New row 37: copy formatting from row 35, then set values
$ws.Range("A35:E35").Copy()
$ws.Range("A37:E37").PasteSpecial(-4122)
$ws.Range("A37").Value = "PublishedAPostTimeStampTest"
$ws.Range("B37").Value = "TBD"
$ws.Range("C37").Value = "Verify that Created Post displayed as per System date"
$ws.Range("D37").Value = "Y"
$ws.Range("E37").Value = "PASS"

# Update selection to match the expanded range
$ws.Range("D2:D37").Select()
